$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Key table (K/L columns): row 7 now documents "door" (type 3), new row 8 documents NPC (type 4) ---
$ws.Range("I4").Value2 = 4
$ws.Range("I5").Value2 = 4
$ws.Range("I6").Value2 = 4

# NOTE: write order matters for new shared-string index assignment -- "Map Doors"
# must become index 68 and "Is a door." index 69, so N11 is written before L7.

# --- Map Construction table (N/O/P columns): insert "Map Doors" row, push "Map data" row down ---
$ws.Range("N11").Value2 = "Map Doors"
$ws.Range("O11").Value2 = 0

$ws.Range("L7").Value2 = "Is a door."

$ws.Range("K8").Value2 = 4
$ws.Range("K8").HorizontalAlignment = -4131
$ws.Range("L8").Value2 = "Is an NPC."

$ws.Range("N12").Value2 = "Map data"
$ws.Range("O12").Value2 = 1
$ws.Range("P12").Value2 = 3600

# --- Door count for the 4th key entry moves from row 13 to row 14 ---
$ws.Range("N13").ClearContents()
$ws.Range("D14").Value2 = 3
$ws.Range("N14").Value2 = "Note: See example below on how the document MUST be formatted."

# --- Move/resize the embedded picture ---
$shp = $ws.Shapes.Item(1)
$shp.Left = 830.560546875
$shp.Top = 225.0
$shp.Width = 241.3603515625
$shp.Height = 201.0

# --- Update the saved selection ---
[void]$ws.Range("G14").Select()
